$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update shared strings used in A6/B6
$ws.Range("A6").Value = "Eintragungsdatum"
$ws.Range("B6").Value = "15.12.2023"

# Update the active cell selection on the sheet
$ws.Activate()
$ws.Range("A3").Select()
